$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing cell values (rows 1-26, columns A-H and J-Q) ---
$ws.Cells.Item(1,1).Value = "negative"
$ws.Cells.Item(1,10).Value = "positive"
$ws.Cells.Item(2,1).Value = "name"
$ws.Cells.Item(2,2).Value = "anchor score"
$ws.Cells.Item(2,3).Value = "type occurences"
$ws.Cells.Item(2,4).Value = "total occurences"
$ws.Cells.Item(2,5).Value = "+%"
$ws.Cells.Item(2,6).Value = "-%"
$ws.Cells.Item(2,7).Value = "both"
$ws.Cells.Item(2,8).Value = "normal"
$ws.Cells.Item(2,10).Value = "name"
$ws.Cells.Item(2,11).Value = "anchor score"
$ws.Cells.Item(2,12).Value = "type occurences"
$ws.Cells.Item(2,13).Value = "total occurences"
$ws.Cells.Item(2,14).Value = "+%"
$ws.Cells.Item(2,15).Value = "-%"
$ws.Cells.Item(2,16).Value = "both"
$ws.Cells.Item(2,17).Value = "normal"
$ws.Cells.Item(3,1).Value = "crude"
$ws.Cells.Item(3,2).Value = 0.8823529411764706
$ws.Cells.Item(3,3).Value = 30
$ws.Cells.Item(3,4).Value = 30
$ws.Cells.Item(3,5).Value = 0
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = $false
$ws.Cells.Item(3,8).Value = 4
$ws.Cells.Item(3,10).Value = "happy"
$ws.Cells.Item(3,11).Value = 1
$ws.Cells.Item(3,12).Value = 26
$ws.Cells.Item(3,13).Value = 26
$ws.Cells.Item(3,14).Value = 1
$ws.Cells.Item(3,15).Value = 0
$ws.Cells.Item(3,16).Value = $false
$ws.Cells.Item(3,17).Value = 0
$ws.Cells.Item(4,1).Value = "crisis"
$ws.Cells.Item(4,2).Value = 0.5856164383561644
$ws.Cells.Item(4,3).Value = 171
$ws.Cells.Item(4,4).Value = 171
$ws.Cells.Item(4,5).Value = 0
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = $false
$ws.Cells.Item(4,8).Value = 121
$ws.Cells.Item(4,10).Value = "interesting"
$ws.Cells.Item(4,11).Value = 0.9393939393939394
$ws.Cells.Item(4,12).Value = 31
$ws.Cells.Item(4,13).Value = 31
$ws.Cells.Item(4,14).Value = 1
$ws.Cells.Item(4,15).Value = 0
$ws.Cells.Item(4,16).Value = $false
$ws.Cells.Item(4,17).Value = 2
$ws.Cells.Item(5,1).Value = "panic"
$ws.Cells.Item(5,2).Value = 0.187984496124031
$ws.Cells.Item(5,3).Value = 97
$ws.Cells.Item(5,4).Value = 97
$ws.Cells.Item(5,5).Value = 0
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = $false
$ws.Cells.Item(5,8).Value = 419
$ws.Cells.Item(5,10).Value = "best"
$ws.Cells.Item(5,11).Value = 0.9152542372881356
$ws.Cells.Item(5,12).Value = 54
$ws.Cells.Item(5,13).Value = 54
$ws.Cells.Item(5,14).Value = 1
$ws.Cells.Item(5,15).Value = 0
$ws.Cells.Item(5,16).Value = $false
$ws.Cells.Item(5,17).Value = 5
$ws.Cells.Item(6,1).Value = "sc"
$ws.Cells.Item(6,2).Value = 0.1534391534391534
$ws.Cells.Item(6,3).Value = 29
$ws.Cells.Item(6,4).Value = 29
$ws.Cells.Item(6,5).Value = 0
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = $false
$ws.Cells.Item(6,8).Value = 160
$ws.Cells.Item(6,10).Value = "great"
$ws.Cells.Item(6,11).Value = 0.8839285714285714
$ws.Cells.Item(6,12).Value = 99
$ws.Cells.Item(6,13).Value = 99
$ws.Cells.Item(6,14).Value = 1
$ws.Cells.Item(6,15).Value = 0
$ws.Cells.Item(6,16).Value = $false
$ws.Cells.Item(6,17).Value = 13
$ws.Cells.Item(7,10).Value = "love"
$ws.Cells.Item(7,11).Value = 0.8695652173913043
$ws.Cells.Item(7,12).Value = 40
$ws.Cells.Item(7,13).Value = 40
$ws.Cells.Item(7,14).Value = 1
$ws.Cells.Item(7,15).Value = 0
$ws.Cells.Item(7,16).Value = $false
$ws.Cells.Item(7,17).Value = 6
$ws.Cells.Item(8,10).Value = "special"
$ws.Cells.Item(8,11).Value = 0.8333333333333334
$ws.Cells.Item(8,12).Value = 30
$ws.Cells.Item(8,13).Value = 30
$ws.Cells.Item(8,14).Value = 1
$ws.Cells.Item(8,15).Value = 0
$ws.Cells.Item(8,16).Value = $false
$ws.Cells.Item(8,17).Value = 6
$ws.Cells.Item(9,10).Value = "thanks"
$ws.Cells.Item(9,11).Value = 0.8048780487804879
$ws.Cells.Item(9,12).Value = 66
$ws.Cells.Item(9,13).Value = 66
$ws.Cells.Item(9,14).Value = 1
$ws.Cells.Item(9,15).Value = 0
$ws.Cells.Item(9,16).Value = $false
$ws.Cells.Item(9,17).Value = 16
$ws.Cells.Item(10,10).Value = "won"
$ws.Cells.Item(10,11).Value = 0.7948717948717948
$ws.Cells.Item(10,12).Value = 31
$ws.Cells.Item(10,13).Value = 31
$ws.Cells.Item(10,14).Value = 1
$ws.Cells.Item(10,15).Value = 0
$ws.Cells.Item(10,16).Value = $false
$ws.Cells.Item(10,17).Value = 8
$ws.Cells.Item(11,10).Value = "thank"
$ws.Cells.Item(11,11).Value = 0.7890625
$ws.Cells.Item(11,12).Value = 101
$ws.Cells.Item(11,13).Value = 101
$ws.Cells.Item(11,14).Value = 1
$ws.Cells.Item(11,15).Value = 0
$ws.Cells.Item(11,16).Value = $false
$ws.Cells.Item(11,17).Value = 27
$ws.Cells.Item(12,10).Value = "positive"
$ws.Cells.Item(12,11).Value = 0.7586206896551724
$ws.Cells.Item(12,12).Value = 44
$ws.Cells.Item(12,13).Value = 44
$ws.Cells.Item(12,14).Value = 1
$ws.Cells.Item(12,15).Value = 0
$ws.Cells.Item(12,16).Value = $false
$ws.Cells.Item(12,17).Value = 14
$ws.Cells.Item(13,10).Value = "free"
$ws.Cells.Item(13,11).Value = 0.7416666666666667
$ws.Cells.Item(13,12).Value = 89
$ws.Cells.Item(13,13).Value = 89
$ws.Cells.Item(13,14).Value = 1
$ws.Cells.Item(13,15).Value = 0
$ws.Cells.Item(13,16).Value = $false
$ws.Cells.Item(13,17).Value = 31
$ws.Cells.Item(14,10).Value = "support"
$ws.Cells.Item(14,11).Value = 0.7264150943396226
$ws.Cells.Item(14,12).Value = 77
$ws.Cells.Item(14,13).Value = 77
$ws.Cells.Item(14,14).Value = 1
$ws.Cells.Item(14,15).Value = 0
$ws.Cells.Item(14,16).Value = $false
$ws.Cells.Item(14,17).Value = 29
$ws.Cells.Item(15,10).Value = "safe"
$ws.Cells.Item(15,11).Value = 0.7253521126760564
$ws.Cells.Item(15,12).Value = 103
$ws.Cells.Item(15,13).Value = 103
$ws.Cells.Item(15,14).Value = 1
$ws.Cells.Item(15,15).Value = 0
$ws.Cells.Item(15,16).Value = $false
$ws.Cells.Item(15,17).Value = 39
$ws.Cells.Item(16,10).Value = "safety"
$ws.Cells.Item(16,11).Value = 0.7058823529411765
$ws.Cells.Item(16,12).Value = 36
$ws.Cells.Item(16,13).Value = 36
$ws.Cells.Item(16,14).Value = 1
$ws.Cells.Item(16,15).Value = 0
$ws.Cells.Item(16,16).Value = $false
$ws.Cells.Item(16,17).Value = 15
$ws.Cells.Item(17,10).Value = "confidence"
$ws.Cells.Item(17,11).Value = 0.6944444444444444
$ws.Cells.Item(17,12).Value = 25
$ws.Cells.Item(17,13).Value = 25
$ws.Cells.Item(17,14).Value = 1
$ws.Cells.Item(17,15).Value = 0
$ws.Cells.Item(17,16).Value = $false
$ws.Cells.Item(17,17).Value = 11
$ws.Cells.Item(18,10).Value = "good"
$ws.Cells.Item(18,11).Value = 0.6875
$ws.Cells.Item(18,12).Value = 110
$ws.Cells.Item(18,13).Value = 110
$ws.Cells.Item(18,14).Value = 1
$ws.Cells.Item(18,15).Value = 0
$ws.Cells.Item(18,16).Value = $false
$ws.Cells.Item(18,17).Value = 50
$ws.Cells.Item(19,10).Value = "heroes"
$ws.Cells.Item(19,11).Value = 0.6595744680851063
$ws.Cells.Item(19,12).Value = 31
$ws.Cells.Item(19,13).Value = 31
$ws.Cells.Item(19,14).Value = 1
$ws.Cells.Item(19,15).Value = 0
$ws.Cells.Item(19,16).Value = $false
$ws.Cells.Item(19,17).Value = 16
$ws.Cells.Item(20,10).Value = "well"
$ws.Cells.Item(20,11).Value = 0.6170212765957447
$ws.Cells.Item(20,12).Value = 58
$ws.Cells.Item(20,13).Value = 58
$ws.Cells.Item(20,14).Value = 1
$ws.Cells.Item(20,15).Value = 0
$ws.Cells.Item(20,16).Value = $false
$ws.Cells.Item(20,17).Value = 36
$ws.Cells.Item(21,10).Value = "fresh"
$ws.Cells.Item(21,11).Value = 0.6041666666666666
$ws.Cells.Item(21,12).Value = 29
$ws.Cells.Item(21,13).Value = 29
$ws.Cells.Item(21,14).Value = 1
$ws.Cells.Item(21,15).Value = 0
$ws.Cells.Item(21,16).Value = $false
$ws.Cells.Item(21,17).Value = 19
$ws.Cells.Item(22,10).Value = "better"
$ws.Cells.Item(22,11).Value = 0.6031746031746031
$ws.Cells.Item(22,12).Value = 38
$ws.Cells.Item(22,13).Value = 38
$ws.Cells.Item(22,14).Value = 1
$ws.Cells.Item(22,15).Value = 0
$ws.Cells.Item(22,16).Value = $false
$ws.Cells.Item(22,17).Value = 25
$ws.Cells.Item(23,10).Value = "relief"
$ws.Cells.Item(23,11).Value = 0.6
$ws.Cells.Item(23,12).Value = 30
$ws.Cells.Item(23,13).Value = 30
$ws.Cells.Item(23,14).Value = 1
$ws.Cells.Item(23,15).Value = 0
$ws.Cells.Item(23,16).Value = $false
$ws.Cells.Item(23,17).Value = 20
$ws.Cells.Item(24,10).Value = "hand"
$ws.Cells.Item(24,11).Value = 0.5326370757180157
$ws.Cells.Item(24,12).Value = 204
$ws.Cells.Item(24,13).Value = 204
$ws.Cells.Item(24,14).Value = 1
$ws.Cells.Item(24,15).Value = 0
$ws.Cells.Item(24,16).Value = $false
$ws.Cells.Item(24,17).Value = 179
$ws.Cells.Item(25,10).Value = "like"
$ws.Cells.Item(25,11).Value = 0.4529411764705882
$ws.Cells.Item(25,12).Value = 154
$ws.Cells.Item(25,13).Value = 154
$ws.Cells.Item(25,14).Value = 1
$ws.Cells.Item(25,15).Value = 0
$ws.Cells.Item(25,16).Value = $false
$ws.Cells.Item(25,17).Value = 186
$ws.Cells.Item(26,10).Value = "care"
$ws.Cells.Item(26,11).Value = 0.4157303370786517
$ws.Cells.Item(26,12).Value = 37
$ws.Cells.Item(26,13).Value = 37
$ws.Cells.Item(26,14).Value = 1
$ws.Cells.Item(26,15).Value = 0
$ws.Cells.Item(26,16).Value = $false
$ws.Cells.Item(26,17).Value = 52

# --- Remove obsolete rows 7-8 from the A:H block (fully deletes the cells) ---
$ws.Range("A7:H8").Clear()

# --- Add new rows 27-32 to the J:Q block ---
$ws.Cells.Item(27,10).Value = "help"
$ws.Cells.Item(27,11).Value = 0.4067796610169492
$ws.Cells.Item(27,12).Value = 120
$ws.Cells.Item(27,13).Value = 120
$ws.Cells.Item(27,14).Value = 1
$ws.Cells.Item(27,15).Value = 0
$ws.Cells.Item(27,16).Value = $false
$ws.Cells.Item(27,17).Value = 175
$ws.Cells.Item(28,10).Value = "protect"
$ws.Cells.Item(28,11).Value = 0.3698630136986301
$ws.Cells.Item(28,12).Value = 27
$ws.Cells.Item(28,13).Value = 27
$ws.Cells.Item(28,14).Value = 1
$ws.Cells.Item(28,15).Value = 0
$ws.Cells.Item(28,16).Value = $false
$ws.Cells.Item(28,17).Value = 46
$ws.Cells.Item(29,10).Value = "please"
$ws.Cells.Item(29,11).Value = 0.3514644351464435
$ws.Cells.Item(29,12).Value = 84
$ws.Cells.Item(29,13).Value = 84
$ws.Cells.Item(29,14).Value = 1
$ws.Cells.Item(29,15).Value = 0
$ws.Cells.Item(29,16).Value = $false
$ws.Cells.Item(29,17).Value = 155
$ws.Cells.Item(30,10).Value = "increase"
$ws.Cells.Item(30,11).Value = 0.3461538461538461
$ws.Cells.Item(30,12).Value = 27
$ws.Cells.Item(30,13).Value = 27
$ws.Cells.Item(30,14).Value = 1
$ws.Cells.Item(30,15).Value = 0
$ws.Cells.Item(30,16).Value = $false
$ws.Cells.Item(30,17).Value = 51
$ws.Cells.Item(31,10).Value = "and"
$ws.Cells.Item(31,11).Value = 0.0101161483701761
$ws.Cells.Item(31,12).Value = 27
$ws.Cells.Item(31,13).Value = 31
$ws.Cells.Item(31,14).Value = 0.87
$ws.Cells.Item(31,15).Value = 0.13
$ws.Cells.Item(31,16).Value = $true
$ws.Cells.Item(31,17).Value = 2642
$ws.Cells.Item(32,10).Value = "the"
$ws.Cells.Item(32,11).Value = 0.004841208365608056
$ws.Cells.Item(32,12).Value = 25
$ws.Cells.Item(32,13).Value = 26
$ws.Cells.Item(32,14).Value = 0.96
$ws.Cells.Item(32,15).Value = 0.04000000000000004
$ws.Cells.Item(32,16).Value = $true
$ws.Cells.Item(32,17).Value = 5139

# --- Apply the existing row-26 formatting to the newly added rows 27-32 ---
$ws.Range("J26:Q26").Copy()
$ws.Range("J27:Q32").PasteSpecial(-4122)
$excel.CutCopyMode = $false
